$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3846153846153846
$ws.Range("C2").Value = 0.3076923076923077
$ws.Range("P2").Value = 0.07692307692307693
$ws.Range("S2").Value = 0.2307692307692308
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.25
$ws.Range("J6").Value = 0.5
$ws.Range("Q6").Value = 0.2
$ws.Range("S6").Value = 0.3
$ws.Range("B7").Value = 0.1818181818181818
$ws.Range("J7").Value = 0.2727272727272727
$ws.Range("R7").Value = 0.09090909090909091
$ws.Range("S7").Value = 0.4545454545454545
$ws.Range("F8").Value = 0.1052631578947368
$ws.Range("J8").Value = 0.1578947368421053
$ws.Range("Q8").Value = 0.4210526315789473
$ws.Range("S8").Value = 0.3157894736842105
$ws.Range("Q9").Value = 0.5
$ws.Range("S9").Value = 0.5
$ws.Range("B10").Value = 0.1129032258064516
$ws.Range("F10").Value = 0.08064516129032258
$ws.Range("J10").Value = 0.1612903225806452
$ws.Range("O10").Value = 0.01612903225806452
$ws.Range("Q10").Value = 0.1935483870967742
$ws.Range("R10").Value = 0.1129032258064516
$ws.Range("S10").Value = 0.3225806451612903
$ws.Range("G11").Value = 0.09523809523809523
$ws.Range("J11").Value = 0.09523809523809523
$ws.Range("K11").Value = 0.2380952380952381
$ws.Range("L11").Value = 0.5238095238095238
$ws.Range("S11").Value = 0.04761904761904762
$ws.Range("G12").Value = 0.6363636363636364
$ws.Range("J12").Value = 0.3636363636363636
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.5
$ws.Range("F15").Value = 0.1666666666666667
$ws.Range("J15").Value = 0.5
$ws.Range("K15").Value = 0.1666666666666667
$ws.Range("S15").Value = 0.1666666666666667
$ws.Range("H16").Value = 0.25
$ws.Range("I16").Value = 0.25
$ws.Range("S16").Value = 0.5
$ws.Range("H17").Value = 0.2173913043478261
$ws.Range("J17").Value = 0.4782608695652174
$ws.Range("K17").Value = 0.04347826086956522
$ws.Range("M17").Value = 0.08695652173913043
$ws.Range("O17").Value = 0.08695652173913043
$ws.Range("S17").Value = 0.08695652173913043
$ws.Range("H18").Value = 0.125
$ws.Range("J18").Value = 0.375
$ws.Range("K18").Value = 0.25
$ws.Range("O18").Value = 0.125
$ws.Range("S18").Value = 0.125
$ws.Range("F19").Value = 0.0196078431372549
$ws.Range("H19").Value = 0.2352941176470588
$ws.Range("I19").Value = 0.0196078431372549
$ws.Range("J19").Value = 0.3137254901960784
$ws.Range("K19").Value = 0.2352941176470588
$ws.Range("M19").Value = 0.0392156862745098
$ws.Range("O19").Value = 0.0196078431372549
$ws.Range("S19").Value = 0.1176470588235294
